# Apply the cryptos-list refresh described in the commit:
# updated prices / % changes, a swapped LidoDAOToken<->TheGraph pair,
# and BitcoinSV replaced by Celestia in the final row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.353.99'
$ws.Range("E2").Value = '  +7.32%  '
$ws.Range("D3").Value = '3.579.51'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '''415.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '''128.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").Value = '''0.648'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.11%  '
$ws.Range("D8").Value = '3.572.41'
$ws.Range("E8").Value = '  +3.12%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '''0.766'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.43%  '
$ws.Range("D11").Value = '''0.175'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +14.76%  '
$ws.Range("E12").Value = '  +45.83%  '
$ws.Range("D13").Value = '''42.03'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '''9.80'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = '4.155.32'
$ws.Range("E15").Value = '  +3.36%  '
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '''20.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '3.575.10'
$ws.Range("E18").Value = '  +3.53%  '
$ws.Range("E19").Value = '  +4.92%  '
$ws.Range("D20").Value = '67.331.70'
$ws.Range("E20").Value = '  +7.35%  '
$ws.Range("D21").Value = '''12.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("D22").Value = '''450.63'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.41%  '
$ws.Range("D23").Value = '''88.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("D24").Value = '''3.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.82%  '
$ws.Range("D25").Value = '''13.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").Value = '''10.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.46%  '
$ws.Range("D28").Value = '''34.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.56%  '
$ws.Range("D29").Value = '''4.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.88%  '
$ws.Range("D30").Value = '''2.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.75%  '
$ws.Range("D31").Value = '''12.28'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").Value = '''7.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.49%  '
$ws.Range("D34").Value = '''0.161'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.60%  '
$ws.Range("D35").Value = '''40.82'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.34%  '
$ws.Range("D36").Value = '''0.998'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '''56.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.82%  '
$ws.Range("D38").Value = '''0.0491'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").Value = '0.0₃0730'
$ws.Range("E39").Value = '  +28.17%  '
$ws.Range("D40").Value = '''0.146'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.34%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '''3.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").Value = '''149.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("B45").Value = 'LidoDAOToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D45").Value = '''3.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '''0.310'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.66%  '
$ws.Range("D47").Value = '''4.28'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.50%  '
$ws.Range("E48").Value = '  -4.55%  '
$ws.Range("D49").Value = '''2.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("D50").Value = '''2.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.78%  '
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = '''15.51'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.27%  '
